$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Clean up the author-list heading (merges runs / removes spell-check
#    artifacts around "Fenoglietto"; visible text is unchanged).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Paul Dubois, Pascal Fenoglietto, Paul-Henry",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paul Dubois, Pascal Fenoglietto, Paul-Henry", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Fix the "ArXiV" contribution heading: drop the leading space and the
#    spell-check split around "ArXiV" (single clean run).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " ArXiV: Radiotherapy Dosimetry: A Review on Open-Source Optimizer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ArXiV: Radiotherapy Dosimetry: A Review on Open-Source Optimizer", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Fix the AIME contribution heading text (drop "(full paper coming
#    soon)" and hyphenate "Knowledge-Based").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "AIME: Radiotherapy Dose Optimization via Clinical Knowledge Based Reinforcement Learning (full paper coming soon)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "AIME: Radiotherapy Dose Optimization via Clinical Knowledge-Based Reinforcement Learning", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Reorder the contributions: move the "SFPM" heading (+ its following
#    blank paragraph) so that it appears after "ASTRO" instead of before
#    "AIME". Final order becomes:
#       ArXiV, ESTRO, AIME, ASTRO, SFPM, SFRO
# ---------------------------------------------------------------------------
function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $paraText = $doc.Paragraphs.Item($i).Range.Text.TrimEnd("`r`a")
        if ($paraText -eq $text) {
            return $i
        }
    }
    return -1
}

$sfpmIdx = Get-ParaIndexByText $d "SFPM: Dose Volume Histograms Guided Deep Dose Predictions"
$sfpmHeading = $d.Paragraphs.Item($sfpmIdx)
$sfpmBlank = $d.Paragraphs.Item($sfpmIdx + 1)

$moveRange = $d.Range($sfpmHeading.Range.Start, $sfpmBlank.Range.End)
$moveRange.Cut() | Out-Null

$astroIdx = Get-ParaIndexByText $d "ASTRO: Clinically Dependent Fully Automatic Treatment Planning System"
$astroBlank = $d.Paragraphs.Item($astroIdx + 1)
$insertPoint = $d.Range($astroBlank.Range.End, $astroBlank.Range.End)
$insertPoint.Paste()

# Cutting/Pasting across paragraph marks can drop paragraph-level style, so
# make sure the moved heading keeps the "Heading 4" style.
$newSfpmIdx = Get-ParaIndexByText $d "SFPM: Dose Volume Histograms Guided Deep Dose Predictions"
$d.Paragraphs.Item($newSfpmIdx).Style = "Heading 4"

# ---------------------------------------------------------------------------
# 5) Fill in the (previously blank) paragraph following each Heading 4 with
#    its one-paragraph summary.
# ---------------------------------------------------------------------------
function Set-SummaryAfterHeading($doc, $headingText, $summaryText) {
    $idx = Get-ParaIndexByText $doc $headingText
    if ($idx -lt 0) {
        throw "Heading not found: $headingText"
    }
    $summaryPara = $doc.Paragraphs.Item($idx + 1)
    $summaryPara.Range.Text = $summaryText
}

Set-SummaryAfterHeading $d `
    "ArXiV: Radiotherapy Dosimetry: A Review on Open-Source Optimizer" `
    "This study evaluates the performance of various state-of-the-art open-source optimizers for radiotherapy dosimetry. Newton CG and LBFGS were the most efficient. These insights help guide the selection of optimization tools for more efficient cancer treatment planning."

Set-SummaryAfterHeading $d `
    "ESTRO: A Novel Framework for Multi-Objective Optimization and Robust Plan Selection Using Graph Theory" `
    "This study presents an innovative framework for optimizing radiotherapy dose distribution by generating and clustering multiple treatment plans with randomized constraint weights. The new proposed framework clusters plans based on dose-volume histogram similarities, which carry most of the clinical meaning."

Set-SummaryAfterHeading $d `
    "AIME: Radiotherapy Dose Optimization via Clinical Knowledge-Based Reinforcement Learning" `
    "This research introduces a deep learning framework for automating radiotherapy treatment planning by training a reinforcement learning agent to mimic dose distributions from past cases. This method represents a first step towards fully automated, human-less treatment planning systems by navigating towards clinically acceptable solutions based on human planners' optimal dose plans."

Set-SummaryAfterHeading $d `
    "ASTRO: Clinically Dependent Fully Automatic Treatment Planning System" `
    "This study demonstrates the potential of training reinforcement learning (RL) agents to mimic human-optimized radiotherapy plans by leveraging past clinical dose data, tailored to specific clinic guidelines. The results suggest that a fully automated treatment planning system (TPS) can be customized for each clinic's practices, improving the feasibility and adoption of automated TPS in clinical settings."

Set-SummaryAfterHeading $d `
    "SFPM: Dose Volume Histograms Guided Deep Dose Predictions" `
    "This study presents a deep-learning model incorporating Dose-Volume Histograms (DVHs) into radiotherapy dose prediction. By integrating target DVH into the model's input, the same model can generate deep doses following a clinical guideline. This technique enables a new workflow where a template of DVHs is used for each clinic, and dosimetrists can fine-tune the target DVHs if needed."

Set-SummaryAfterHeading $d `
    "SFRO: Attention Mechanism on Dose-Volume Histograms for Deep Dose Predictions" `
    "This study introduces a new approach for radiotherapy dose prediction by incorporating Dose-Volume Histograms (DVHs) into deep learning models using an attention mechanism. This approach slightly improves dose prediction accuracy."

# ---------------------------------------------------------------------------
# 6) Touch the built-in "No Spacing" style so it gets registered in the
#    document (then revert back to Normal so no paragraph references it).
# ---------------------------------------------------------------------------
$lastIdx = Get-ParaIndexByText $d "SFRO: Attention Mechanism on Dose-Volume Histograms for Deep Dose Predictions"
$touchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$touchPara.Style = "No Spacing"
$touchPara.Style = "Normal"
